$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the current row 31 (before row 32).
# Excel shifts rows 32:43 down to 34:45, preserving their content/formatting,
# which also reproduces the duplicated tail rows (new 44:45 == old 42:43).
$ws.Rows("32:33").Insert()

# --- New row 32 ---
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "Femacal de La Calera"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = "2023-08-23"
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = 100112043
$ws.Range("G32").Value = "Pepino dulce"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 45
$ws.Range("K32").Value = 23000
$ws.Range("L32").Value = 23000
$ws.Range("M32").Value = 23000
$ws.Range("N32").Value = "`$/caja 15 kilos"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 1533
$ws.Range("Q32").Value = 15
$ws.Range("R32").Value = "Hortaliza"

# --- New row 33 ---
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = "Femacal de La Calera"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = "2023-08-23"
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 100112043
$ws.Range("G33").Value = "Pepino dulce"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 17000
$ws.Range("N33").Value = "`$/caja 15 kilos"
$ws.Range("O33").Value = "Provincia de Limarí"
$ws.Range("P33").Value = 1133
$ws.Range("Q33").Value = 15
$ws.Range("R33").Value = "Hortaliza"
